$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.093.52"
$ws.Range("E2").Value = "  -0.86%  "

$ws.Range("D3").Value = "2.316.99"
$ws.Range("E3").Value = "  -1.82%  "

$ws.Range("E4").Value = "  -0.09%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "533.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.56"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.98%  "

$ws.Range("E7").Value = "  -0.27%  "

$ws.Range("E8").Value = "  -0.64%  "

$ws.Range("D9").Value = "2.338.39"
$ws.Range("E9").Value = "  -1.68%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.102"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.40%  "

$ws.Range("E11").Value = "  +0.00%  "

$ws.Range("E12").Value = "  -2.74%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.347"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.93%  "

$ws.Range("D14").Value = "2.730.28"
$ws.Range("E14").Value = "  -1.86%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "23.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -3.73%  "

$ws.Range("D16").Value = "57.100.93"
$ws.Range("E16").Value = "  -0.84%  "

$ws.Range("E17").Value = "  -2.64%  "

$ws.Range("D18").Value = "2.323.21"
$ws.Range("E18").Value = "  -1.81%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "337.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.95%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "10.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.03%  "

$ws.Range("E21").Value = "  +2.68%  "

$ws.Range("E22").Value = "  -2.30%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.17%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "61.67"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.167"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.60%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.68"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +4.30%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.994"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.15%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.33"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.78%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.62%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.73"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.06%  "

$ws.Range("D31").Value = "0.0₃0726"
$ws.Range("E31").Value = "  -3.01%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.09"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.29%  "

$ws.Range("E33").Value = "  -0.55%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.991"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.24%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.26"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.93%  "

$ws.Range("B37").Value = "SuiNetwork"
$ws.Range("C37").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.907"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.47%  "

$ws.Range("B38").Value = "NEARProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.41%  "

$ws.Range("E39").Value = "  -0.20%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "39.17"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.26%  "

$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.77"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +8.05%  "

$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "148.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.74%  "

$ws.Range("B43").Value = "PolygonEcosystemToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.375"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.44%  "

$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.60"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.83%  "

$ws.Range("B45").Value = "Bittensor"
$ws.Range("C45").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "281.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.89%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0927"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.51%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0502"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.45%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "18.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.557"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0217"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.24%  "

$ws.Range("E51").Value = "  +7.82%  "
